# Commit: "finalized two designs (resume - bottom and resume - back)"
# - Insert an en dash run immediately before the "EDUCATION" heading run
#   (the very first run in the document), then place the "_GoBack"
#   bookmark right after that dash / right before "EDUCATION".
# - A document only ever has a single "_GoBack" bookmark, so re-adding it
#   here automatically retires the old one that used to sit in the empty
#   paragraph just above "MORE INFORMATION" - which is exactly the other
#   half of this diff.
$d = $word.ActiveDocument
$eduStart = $d.Range(0, 0)
$eduStart.InsertBefore([char]0x2013)
$d.Bookmarks.Add("_GoBack", $d.Range(1, 1))
